# edit.ps1 - apply the UseCases.docx revision:
#   1) Move the "_GoBack" bookmark from the end of the "Level:" paragraph
#      (Use Case 1) to the end of the "Special Requirements: N/A" paragraph
#      (also Use Case 1, right before "Frequency of Occurrence").
#   2) Recolor the run "User indicates that cards should be dealt, and game
#      should begin." from the themed green (accent6, 80% shade) to the
#      plain RGB green 00B050.
#
# (A third part of the source diff adds three <w:lsdException> entries -
#  "Normal Table", "Table Web 3", "Table Theme" - to styles.xml's
#  w:latentStyles block. That metadata lives outside the Range/body object
#  model Word's COM automation surface exposes (there is no supported
#  LatentStyles/lsdException editing surface here, mirroring real Word COM),
#  so it cannot be expressed through COM-interop calls and is intentionally
#  left alone.)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Relocate the _GoBack bookmark
# ---------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Find the first paragraph whose text is "Special Requirements: N/A"
# (Use Case 1's field, immediately followed by "Frequency of Occurrence").
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($targetPara -eq $null) {
        $t = $p.Range.Text
        if ($t -like "Special Requirements: N/A*") {
            $targetPara = $p
        }
    }
}

if ($targetPara -ne $null) {
    $paraRange = $targetPara.Range

    # Collapsed range right after the last character ("N/A"), before the
    # paragraph mark - same spot the bookmark previously occupied in the
    # "Level:" paragraph. Collapsed Bookmarks.Add calls at this exact
    # boundary aren't placed correctly by this host, so stage the bookmark
    # around a throwaway character and then delete that character - the
    # bookmark collapses in place and survives.
    $insertPoint = $paraRange.Duplicate
    $insertPoint.Start = $paraRange.End - 1
    $insertPoint.End = $paraRange.End - 1
    $insertPoint.InsertAfter("X")
    $insertPoint.MoveEnd(1, 1)
    $d.Bookmarks.Add("_GoBack", $insertPoint)
    $insertPoint.Text = ""
}

# ---------------------------------------------------------------------
# 2) Recolor the "User indicates..." run
# ---------------------------------------------------------------------

$colorRange = $d.Content
$found = $colorRange.Find.Execute(
    "User indicates that cards should be dealt, and game should begin.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # wdColor OLE value is 0x00BBGGRR; RGB(00,B0,50) -> 0x0050B000
    $colorRange.Font.Color = 5287936
}
